$wb = $excel.ActiveWorkbook

# --- Update conversion text on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.18 = 12172.93 pesos`n✅ 12172.93 pesos = 3.16 = 968.96 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update rate figures on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 314
$wsTasas.Range("O10").Value = 3822.3
$wsTasas.Range("N12").Value = 3856.79
$wsTasas.Range("O12").Value = 307
